$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text (string) type rather than being
# auto-converted to numbers/dates by Excel, matching the original inlineStr cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.329.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.271.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '499.05'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.39'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0956'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.337'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.93'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.671.64'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.267.77'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.265.55'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.29'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.88%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '304.47'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '60.59'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.37'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '175.33'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0706'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.01'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.12%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.61'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.45%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.85'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.951'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.48%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.73'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.41'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.84'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0898'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '245.88'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.25'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.65%  '
